$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the task name / abbreviation on row 2 to use the new response codes
$ws.Range("B2").Value = "NEWTD1"
$ws.Range("A2").Value = "New Task Def 1"

# Add a new column header for tutorial_stream
$ws.Range("S1").Value = "tutorial_stream"

# Update the active selection to match the new last column
$ws.Range("S1").Select()
